$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (index 1 / rId1) - rows 6-40
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 380
$ws1.Range("F7").Value = 1203
$ws1.Range("F8").Value = 452
$ws1.Range("F9").Value = 7392
$ws1.Range("F10").Value = 85
$ws1.Range("F11").Value = 94
$ws1.Range("F12").Value = 2060
$ws1.Range("F13").Value = 8049
$ws1.Range("F16").Value = 5535
$ws1.Range("F17").Value = 56
$ws1.Range("F18").Value = 2457
$ws1.Range("F19").Value = 1042
$ws1.Range("F21").Value = 308
$ws1.Range("F22").Value = 388
$ws1.Range("F24").Value = 14
$ws1.Range("F25").Value = 404
$ws1.Range("F26").Value = 313
$ws1.Range("F27").Value = 15
$ws1.Range("F28").Value = 2474
$ws1.Range("F31").Value = 91
$ws1.Range("F32").Value = 172
$ws1.Range("F33").Value = 607
$ws1.Range("F36").Value = 1544
$ws1.Range("F38").Value = 9
$ws1.Range("F39").Value = 2422
$ws1.Range("F40").Value = 2230

# Sheet 2: "演出" (index 2 / rId2) - rows 2-5
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 106
$ws2.Range("F3").Value = 86
$ws2.Range("F4").Value = 80
$ws2.Range("F5").Value = 21

# Sheet 4: "全部类型" (index 4 / rId4) - rows 6-45
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value = 106
$ws4.Range("F7").Value = 380
$ws4.Range("F8").Value = 1203
$ws4.Range("F9").Value = 452
$ws4.Range("F10").Value = 7392
$ws4.Range("F11").Value = 85
$ws4.Range("F12").Value = 94
$ws4.Range("F13").Value = 2060
$ws4.Range("F14").Value = 8049
$ws4.Range("F17").Value = 5535
$ws4.Range("F18").Value = 56
$ws4.Range("F19").Value = 2457
$ws4.Range("F20").Value = 1042
$ws4.Range("F22").Value = 388
$ws4.Range("F24").Value = 86
$ws4.Range("F25").Value = 14
$ws4.Range("F26").Value = 80
$ws4.Range("F27").Value = 404
$ws4.Range("F28").Value = 313
$ws4.Range("F29").Value = 15
$ws4.Range("F30").Value = 2474
$ws4.Range("F33").Value = 91
$ws4.Range("F34").Value = 172
$ws4.Range("F35").Value = 21
$ws4.Range("F36").Value = 607
$ws4.Range("F40").Value = 1544
$ws4.Range("F42").Value = 9
$ws4.Range("F43").Value = 2422
$ws4.Range("F45").Value = 2230
